$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are treated as plain text (matches source inlineStr "t" cells)
# so Excel does not auto-convert numeric-looking / percent-looking strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.28%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.46%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.160"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.18%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05802"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.59%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.668"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.60%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.228"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "6.96%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.62%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8604"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.53%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1376"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.04%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07089"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.92%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03206"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "10.68%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09380"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.25%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001531"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.20%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005909"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.10%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.496"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.49%"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.226"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.15%"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006020"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.22%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3196"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.59%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03354"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.57%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.71%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.07%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04133"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.99%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001228"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.77%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004140"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.51%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "4.10%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03752"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.34%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005690"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "64.60%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1071"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.31%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002200"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.19%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009578"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.04%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005291"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.39%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.22%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05800"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-35.42%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-20.93%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.22%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.22%"
